$d = $word.ActiveDocument

# --- First paragraph: the "**ID__AFFARS_..." placeholder paragraph ---
$p = $d.Paragraphs.Item(1)

# Add a paragraph border (top/left/bottom/right) with 5pt space, no line
$borders = $p.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# Replace the placeholder text (including the trailing run's space) with
# the updated topic id, collapsing the two runs into a single run.
$rng = $p.Range
$rng.Find.Execute("**ID__AFFARS_5350_topic_3__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5350_101_1__ID**", 2)
